$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name and title to reflect new "through" date (10-29 -> 10-30)
$ws.Name = "Through 2022-10-30"

# Update the header label for the October row
$ws.Range("A11").Value = "October (through 10-30)"

# September row (row 10) - only the 2022 column (I) changes
$ws.Range("I10").Value = 144

# October row (row 11) - all year columns change
$ws.Range("B11").Value = 30
$ws.Range("C11").Value = 56
$ws.Range("D11").Value = 79
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = 59
$ws.Range("G11").Value = 148
$ws.Range("H11").Value = 186
$ws.Range("I11").Value = 121

# Total row (row 12) - all year columns change
$ws.Range("B12").Value = 256
$ws.Range("C12").Value = 485
$ws.Range("D12").Value = 706
$ws.Range("E12").Value = 613
$ws.Range("F12").Value = 481
$ws.Range("G12").Value = 1049
$ws.Range("H12").Value = 1433
$ws.Range("I12").Value = 1397
